# Add a "Subclass" column to the Players sheet (right after "Class"),
# populate it with each player's subclass, and clear out the old
# Subclass values that used to live in column P (now shifted to Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Insert a new column before the current column D ("Player"), shifting
# D..Q right to E..R. The new column inherits formatting from column C.
$ws.Columns("D:D").Insert()

# New column D header + subclass values.
$ws.Range("D1").Value = "Subclass"
$ws.Range("D2").Value = "Assassin"
$ws.Range("D3").Value = "Beast Master"
$ws.Range("D4").Value = "Lore"
$ws.Range("D5").Value = "Hexblade"
$ws.Range("D6").Value = "Vengeance"

# The old Subclass column (previously P) is now Q after the insert -
# clear its contents, leaving an empty column between Initiative (P)
# and Alignment (R).
$ws.Range("Q1:Q6").ClearContents()

# Update the active selection as recorded in the edited workbook.
$ws.Range("D11").Select()
